$wb = $excel.ActiveWorkbook

# The battery current labels live on the "Add Panels" sheet (J8/K8)
$ws = $wb.Worksheets.Item("Add Panels")

$ws.Range("J8").Value = "Alarm Current(A)"
$ws.Range("K8").Value = "Standby Current(A)"

# Reflect the author's final cell selection on this sheet
$ws.Activate()
$ws.Range("B7").Select()
